$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the A73 value with the precise floating point serial used by the source system
$ws.Range("A73").Value = 44386.76955109606

# Append new row 74
$ws.Range("A74").Value = 44387.76718575539
$ws.Range("B74").Value = 80150
$ws.Range("C74").Value = 67671
$ws.Range("D74").Value = 3606
$ws.Range("E74").Value = 2211
$ws.Range("F74").Value = 1576
$ws.Range("G74").Value = 21291
$ws.Range("H74").Value = 1556
$ws.Range("I74").Value = 895
$ws.Range("J74").Value = 199

# Match style of column A (date format) from row 73
$ws.Range("A74").NumberFormat = $ws.Range("A73").NumberFormat
